# Atualiza modelo e graficos para projeto escrito
# Remove the "NO_MUNICIPIO_RESIDENCIA" row from the tipo_colunas table.
# This was row 7 (A7=5, B7="NO_MUNICIPIO_RESIDENCIA",
# C7="VARIAVEL QUALITATIVA NOMINAL POLICOTOMICA"). Deleting the whole row
# shifts every row below it up by one and Excel keeps the sequential
# index in column A contiguous (0..18 instead of 0..19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Delete()

# Column A holds a contiguous 0-based index; after removing the row the
# values below the deletion point must be renumbered so there is no gap.
for ($r = 7; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
